# Rebuild the "all_exp_model_stats" table: the CEWL model is unchanged, but the
# body condition / osmolality / hematocrit models were refit without the 2-day data,
# splitting the old "day" predictor rows into separate "day (mid)" and "day (after)"
# rows, and a blank spacer row is added before the hematocrit section.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate over the old A1:G24 range, then rebuild through row 28.
$ws.Range("A1:G28").Clear()

# Column headers (row 1) - General format, right aligned (unchanged from before)
$ws.Range("C1").Value = "estimate"
$ws.Range("D1").Value = "SE"
$ws.Range("E1").Value = "t-value"
$ws.Range("F1").Value = "df"
$ws.Range("G1").Value = "p-value"
$ws.Range("C1:G1").HorizontalAlignment = -4152

# "CEWL ~" section title (row 2) - unaffected by this edit
$ws.Range("A2").Value = "CEWL ~"

# Data rows: Row = worksheet row, B = row label, C/D/E/F = estimate/SE/t-value/df,
# G = p-value (either text like "< 0.0001"/"< 0.001", or a precise numeric p-value)
$dataRows = @(
    @{ Row=3; B="(intercept)"; C=-69.375784903449699; D=11.8197221990293; E=-5.8694936932736796; F=230.96964219851; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=4; B="day (after)"; C=18.263176469806002; D=1.9250825011117001; E=9.4869578105142995; F=304.95405323256699; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=5; B="humidity treatment (dry)"; C=0.76053159097103895; D=2.7927378035519199; E=0.27232473811317398; F=41.133118618782198; G=0.78673473866406696; GIsText=$false; GFormat="0.00" }
    @{ Row=6; B="region (ventrum)"; C=15.7650562809789; D=1.9395728079888199; E=8.1281074966842493; F=286.009477668572; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=7; B="region (head)"; C=2.9137684021909802; D=1.9395728079888199; E=1.50227327903835; F=286.00947767441397; G=0.13412980799506699; GIsText=$false; GFormat="0.00" }
    @{ Row=8; B="region (dewlap)"; C=3.0976537768641998; D=1.9482735061366401; E=1.5899481089832901; F=286.09695220181101; G=0.11295088865526801; GIsText=$false; GFormat="0.00" }
    @{ Row=9; B="region (mite patch)"; C=3.7911890489553799; D=1.95647605044726; E=1.9377640978987201; F=286.27000179867201; G=0.053635293484382902; GIsText=$false; GFormat="0.00" }
    @{ Row=10; B="cloacal temperature at the time of measurement (℃)"; C=3.8920479697064101; D=0.48858126829931797; E=7.9660196209610703; F=279.63437114745; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=11; B="day (after) * humidity treatment (dry)"; C=-20.1375496328845; D=2.48807741213446; E=-8.0936186047398806; F=288.34291234635401; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=13; B="(intercept)"; C=10.3769451760934; D=0.22439268559085801; E=46.244578555532897; F=100.999999997624; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=14; B="day (mid)"; C=-0.61728688527405395; D=0.27282504584974299; E=-2.26257410990786; F=100.99999999777199; G=0.025803006571743298; GIsText=$false; GFormat="0.00" }
    @{ Row=15; B="day (after)"; C=-1.0308219364544899; D=0.27282504584974299; E=-3.7783258983569001; F=100.99999999777199; G="< 0.001"; GIsText=$true; GFormat="0.000" }
    @{ Row=16; B="humidity treatment (dry)"; C=0.45433159623644898; D=0.22285169557433099; E=2.0387172512443801; F=100.999999998102; G=0.0440909258449107; GIsText=$false; GFormat="0.00" }
    @{ Row=18; B="(intercept)"; C=368.294322927167; D=10.094991376666799; E=36.4828764270598; F=4.4846729060523796; G=0.00000101210420087086; GIsText=$false; GFormat="0.00" }
    @{ Row=19; B="day (mid)"; C=-17.176470588234899; D=6.8064122486609104; E=-2.5235718849698201; F=92.951110759141301; G=0.0133142302023167; GIsText=$false; GFormat="0.00" }
    @{ Row=20; B="day (after)"; C=3.8980579089779601; D=7.03862810875831; E=0.55380932885593503; F=92.971889700645505; G=0.58103887555425104; GIsText=$false; GFormat="0.00" }
    @{ Row=21; B="humidity treatment (dry)"; C=10.2149493910569; D=6.7131580975050102; E=1.52163098837987; F=92.955743227385597; G=0.13149614644426999; GIsText=$false; GFormat="0.00" }
    @{ Row=22; B="day (mid) * humidity treatment (dry)"; C=-12.3235294117653; D=9.4910884209441502; E=-1.2984316303039301; F=92.951110761998905; G=0.19735132595211199; GIsText=$false; GFormat="0.00" }
    @{ Row=23; B="day (after) * humidity treatment (dry)"; C=-17.783329455543502; D=9.7207971543447496; E=-1.82941061038344; F=92.954413140944595; G=0.070543972502939006; GIsText=$false; GFormat="0.00" }
    @{ Row=26; B="(intercept)"; C=35.591574255057402; D=1.7116680193613401; E=20.793503093162499; F=4.8534561623837096; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=27; B="day (mid)"; C=-5.7714285714285696; D=1.4722666473739801; E=-3.9200973422327001; F=98.766491740676102; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
    @{ Row=28; B="day (after)"; C=-8.0571428571428498; D=1.4722666473739801; E=-5.4726111411367402; F=98.766491734250806; G="< 0.0001"; GIsText=$true; GFormat="0.00" }
)

foreach ($r in $dataRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    # estimate/SE/t-value/df: "0.00", general (left-default) alignment
    $ws.Range("C" + $r.Row + ":F" + $r.Row).NumberFormat = "0.00"
    $gCell = $ws.Cells.Item($r.Row, 7)
    $gCell.Value = $r.G
    if ($r.GIsText) {
        # "< 0.0001" / "< 0.001" text: right aligned, matching the old threshold cells
        $gCell.HorizontalAlignment = -4152
    } else {
        # precise numeric p-value: "0.00"/"0.000", general alignment
        $gCell.NumberFormat = $r.GFormat
    }
}

# Section title rows: "<model> ~" label in column A, blank right-aligned "0.00" cells C:F
$sectionRows = @(
    @{ Row=12; A="body condition ~" }
    @{ Row=17; A="osmolality ~" }
    @{ Row=25; A="hematocrit ~" }
)
foreach ($r in $sectionRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $blankRange = $ws.Range("C" + $r.Row + ":F" + $r.Row)
    $blankRange.Value = ""
    $blankRange.NumberFormat = "0.00"
    $blankRange.HorizontalAlignment = -4152
}

# Fully blank separator row between the osmolality and hematocrit sections (row 24)
$blank24 = $ws.Range("C24:F24")
$blank24.Value = ""
$blank24.NumberFormat = "0.00"
$blank24.HorizontalAlignment = -4152
$ws.Range("G24").Value = ""
$ws.Range("G24").NumberFormat = "0.000"
$ws.Range("G24").HorizontalAlignment = -4152

# Match the author's final selection/cursor position
$ws.Range("I26").Select()
